# Weekly fruit/vegetable update: insert a new weekly price record as row 38,
# pushing the existing rows 38-68 down to rows 39-69.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 38 (shifts rows 38..68 down to 39..69)
$ws.Rows.Item(38).Insert()

# Populate the newly inserted row 38 with the new weekly record
$ws.Range("A38").Value = 3
$ws.Range("B38").Value = "Femacal de La Calera"
$ws.Range("C38").Value = "Coquimbo"
$ws.Range("D38").Value = 44874
$ws.Range("E38").Value = 5
$ws.Range("F38").Value = 100112022
$ws.Range("G38").Value = "Arveja Verde"
$ws.Range("H38").Value = "Perfection"
$ws.Range("I38").Value = "Primera"
$ws.Range("J38").Value = 78
$ws.Range("K38").Value = 20000
$ws.Range("L38").Value = 21000
$ws.Range("M38").Value = 20513
$ws.Range("N38").Value = '$/malla 25 kilos'
$ws.Range("O38").Value = 'Provincia de Limarí'
$ws.Range("P38").Value = 821
$ws.Range("Q38").Value = 25
$ws.Range("R38").Value = "Hortaliza"
